# Supplemental table 2 (2018 data) - update literature RI values (col M)
# and mark their Source (col N) as NIST where applicable.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compound annotation updated")

$ws.Range("M11").Value = 853
$ws.Range("N11").Value = "NIST"
$ws.Range("M12").Value = 926
$ws.Range("N12").Value = "NIST"
$ws.Range("M17").Value = 989
$ws.Range("N17").Value = "NIST"
$ws.Range("M50").Value = 1310
$ws.Range("N50").Value = "NIST"
$ws.Range("M54").Value = 1023
$ws.Range("N54").Value = "NIST"
$ws.Range("M75").Value = 1510
$ws.Range("N75").Value = "NIST"
$ws.Range("M84").Value = 946
$ws.Range("N84").Value = "NIST"
$ws.Range("M101").Value = 1182
$ws.Range("N101").Value = "NIST"
$ws.Range("M103").Value = 1021
$ws.Range("N103").Value = "NIST"
$ws.Range("M106").Value = 1436
$ws.Range("N106").Value = "NIST"
$ws.Range("M136").Value = 758
$ws.Range("N136").Value = "NIST"
$ws.Range("M137").Value = 883
$ws.Range("N137").Value = "NIST"
$ws.Range("M140").Value = 1045
$ws.Range("N140").Value = "NIST"
$ws.Range("M142").Value = 793
$ws.Range("N142").Value = "NIST"
$ws.Range("M148").Value = 843
$ws.Range("N148").Value = "NIST"
$ws.Range("M150").Value = 966
$ws.Range("N150").Value = "NIST"
$ws.Range("M157").Value = 815
$ws.Range("N157").Value = "NIST"
$ws.Range("M159").Value = 981
$ws.Range("N159").Value = "NIST"
$ws.Range("M162").Value = 1774
$ws.Range("N162").Value = "NIST"
$ws.Range("M165").Value = 1492
$ws.Range("N165").Value = "NIST"
$ws.Range("M194").Value = 1308
$ws.Range("N194").Value = "NIST"
$ws.Range("M197").Value = 960
$ws.Range("N197").Value = "NIST"
$ws.Range("M198").Value = 1121
$ws.Range("N198").Value = "NIST"
$ws.Range("M200").Value = 1080
$ws.Range("N200").Value = "NIST"
$ws.Range("M203").Value = 1094
$ws.Range("N203").Value = "NIST"
$ws.Range("M205").Value = 1735
$ws.Range("N205").Value = "NIST"
$ws.Range("M214").Value = 1377
$ws.Range("M216").Value = 909
$ws.Range("N216").Value = "NIST"
$ws.Range("M226").Value = 1531
$ws.Range("N226").Value = "NIST"
$ws.Range("M228").Value = 1253
$ws.Range("N228").Value = "NIST"
$ws.Range("M237").Value = 923
$ws.Range("M240").Value = 843
$ws.Range("N240").Value = "NIST"
$ws.Range("M245").Value = 1060
$ws.Range("N245").Value = "NIST"
$ws.Range("M255").Value = 1159
$ws.Range("N255").Value = "NIST"
$ws.Range("M259").Value = 1096
$ws.Range("N259").Value = "NIST"
$ws.Range("M261").Value = 921
$ws.Range("N261").Value = "NIST"
$ws.Range("M266").Value = 959
$ws.Range("N266").Value = "NIST"
$ws.Range("M268").Value = 1355
$ws.Range("N268").Value = "NIST"
$ws.Range("M270").Value = 1302
$ws.Range("N270").Value = "NIST"
$ws.Range("M273").Value = 1172
$ws.Range("N273").Value = "NIST"
$ws.Range("M275").Value = 777
$ws.Range("N275").Value = "NIST"
$ws.Range("M278").Value = 856
$ws.Range("N278").Value = "NIST"
$ws.Range("M279").Value = 1583
$ws.Range("N279").Value = "NIST"
$ws.Range("M281").Value = 1045
$ws.Range("M293").Value = 1167
$ws.Range("N293").Value = "NIST"
$ws.Range("M298").Value = 1655
$ws.Range("N298").Value = "NIST"
$ws.Range("M299").Value = 710
$ws.Range("N299").Value = "NIST"
$ws.Range("M305").Value = 1401
$ws.Range("N305").Value = "NIST"
$ws.Range("M309").Value = 866
$ws.Range("N309").Value = "NIST"
$ws.Range("M312").Value = 758
$ws.Range("N312").Value = "NIST"
$ws.Range("M324").Value = 1233
$ws.Range("N324").Value = "NIST"

# Restore the author's on-save cursor/scroll position.
$ws.Range("M163").Select()
$excel.ActiveWindow.ScrollRow = 148
$excel.ActiveWindow.ScrollColumn = 1
